$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

# Numeric columns
$ws.Cells.Item($row, 1).Value = 112063112          # A  Id
$ws.Cells.Item($row, 2).Value = 57580              # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value = 208249             # E  TaxonId
$ws.Cells.Item($row, 17).Value = 580550.1418178778 # Q  Ost
$ws.Cells.Item($row, 18).Value = 6579521.799528075 # R  Nord
$ws.Cells.Item($row, 19).Value = 10                # S  Noggrannhet

# Text columns
$ws.Cells.Item($row, 3).Value = "Ovaliderad"          # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value = "LC"                  # D  Rödlistade
$ws.Cells.Item($row, 6).Value = "Vanlig groda"        # F  Artnamn
$ws.Cells.Item($row, 7).Value = "Rana temporaria"     # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "Linnaeus, 1758"      # H  Auktor

# I (Antal) is stored as text "1" in the source data, not a number
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "1"                   # I  Antal

# Blank text placeholders (present in source as empty cells)
$ws.Cells.Item($row, 10).NumberFormat = "@"
$ws.Cells.Item($row, 10).Value = ""    # J  Enhet
$ws.Cells.Item($row, 11).NumberFormat = "@"
$ws.Cells.Item($row, 11).Value = ""    # K  Ålder-Stadium
$ws.Cells.Item($row, 12).NumberFormat = "@"
$ws.Cells.Item($row, 12).Value = ""    # L  Kön
$ws.Cells.Item($row, 13).NumberFormat = "@"
$ws.Cells.Item($row, 13).Value = ""    # M  Aktivitet
$ws.Cells.Item($row, 14).NumberFormat = "@"
$ws.Cells.Item($row, 14).Value = ""    # N  Metod

$ws.Cells.Item($row, 16).Value = "Flugmötesskogen, Hagalund, Srm"  # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Södermanland"       # T  Län
$ws.Cells.Item($row, 21).Value = "Eskilstuna"         # U  Kommun
$ws.Cells.Item($row, 22).Value = "Södermanland"       # V  Provins
$ws.Cells.Item($row, 23).Value = "Eskilstuna"         # W  Församling

# Dates/times are stored as plain text, not real date/time values
$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2023-09-10"         # Y  Startdatum
$ws.Cells.Item($row, 26).NumberFormat = "@"
$ws.Cells.Item($row, 26).Value = "11:00"              # Z  Starttid
$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2023-09-10"         # AA Slutdatum
$ws.Cells.Item($row, 28).NumberFormat = "@"
$ws.Cells.Item($row, 28).Value = "11:00"              # AB Sluttid

# Booleans
$ws.Cells.Item($row, 30).Value = $false  # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false  # AE Osäker artbestämning

$ws.Cells.Item($row, 32).NumberFormat = "@"
$ws.Cells.Item($row, 32).Value = ""      # AF Bestämningsmetod (blank)

$ws.Cells.Item($row, 33).Value = $false  # AG Ospontan

$ws.Cells.Item($row, 34).Value = "Skogsmark"          # AH Biotop
$ws.Cells.Item($row, 35).Value = "Kontinuitetsskog"   # AI Biotop-beskrivning

$ws.Cells.Item($row, 46).NumberFormat = "@"
$ws.Cells.Item($row, 46).Value = ""      # AT Bestämningsår (blank)

$ws.Cells.Item($row, 49).Value = "Michael Lander"  # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Michael Lander"  # AX Observatörer

$ws.Cells.Item($row, 51).NumberFormat = "@"
$ws.Cells.Item($row, 51).Value = ""      # AY Projektnamn (blank)
